$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.811.27'
$ws.Range('E2').Value = '  +3.95%  '
$ws.Range('D3').Value = '1.913.83'
$ws.Range('E3').Value = '  +1.44%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '''250.27'
$ws.Range('E5').Value = '  +1.41%  '
$ws.Range('D6').Value = '''0.701'
$ws.Range('E6').Value = '  +0.21%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +8.16%  '
$ws.Range('D9').Value = '''0.374'
$ws.Range('E9').Value = '  +4.70%  '
$ws.Range('D10').Value = '''58.07'
$ws.Range('E10').Value = '  +8.41%  '
$ws.Range('D11').Value = '''0.0761'
$ws.Range('E11').Value = '  +1.48%  '
$ws.Range('E12').Value = '  +2.02%  '
$ws.Range('E13').Value = '  +8.12%  '
$ws.Range('D14').Value = '''0.817'
$ws.Range('E14').Value = '  +5.65%  '
$ws.Range('D15').Value = '2.192.21'
$ws.Range('D16').Value = '''5.12'
$ws.Range('E16').Value = '  +3.16%  '
$ws.Range('D17').Value = '1.916.60'
$ws.Range('E17').Value = '  +1.57%  '
$ws.Range('D18').Value = '37.271.77'
$ws.Range('E18').Value = '  +5.32%  '
$ws.Range('D19').Value = '''74.69'
$ws.Range('E19').Value = '  +1.44%  '
$ws.Range('D20').Value = '0.0₃0858'
$ws.Range('E20').Value = '  +3.49%  '
$ws.Range('D21').Value = '''13.62'
$ws.Range('E21').Value = '  +5.73%  '
$ws.Range('D22').Value = '''251.34'
$ws.Range('E22').Value = '  +2.34%  '
$ws.Range('D23').Value = '''5.18'
$ws.Range('E23').Value = '  -0.75%  '
$ws.Range('E24').Value = '  -1.02%  '
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('E26').Value = '  +3.17%  '
$ws.Range('D27').Value = '''167.86'
$ws.Range('E27').Value = '  +1.63%  '
$ws.Range('D28').Value = '''8.81'
$ws.Range('E28').Value = '  +1.48%  '
$ws.Range('E29').Value = '  +1.96%  '
$ws.Range('D30').Value = '''0.129'
$ws.Range('E30').Value = '  +0.24%  '
$ws.Range('D31').Value = '''4.60'
$ws.Range('E31').Value = '  +6.88%  '
$ws.Range('D32').Value = '''0.0619'
$ws.Range('E32').Value = '  +3.68%  '
$ws.Range('D33').Value = '''4.33'
$ws.Range('E33').Value = '  +3.16%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').Value = '''1.91'
$ws.Range('E34').Value = '  +1.71%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = '''0.0884'
$ws.Range('E35').Value = '  +19.68%  '
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('D37').Value = '''19.55'
$ws.Range('E37').Value = '  +60.39%  '
$ws.Range('D38').Value = '''1.50'
$ws.Range('E38').Value = '  +1.21%  '
$ws.Range('D39').Value = '''0.883'
$ws.Range('E39').Value = '  +2.84%  '
$ws.Range('E40').Value = '  +2.41%  '
$ws.Range('D41').Value = '''105.47'
$ws.Range('E41').Value = '  +7.94%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = '''0.0229'
$ws.Range('E42').Value = '  +4.67%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').Value = '''18.00'
$ws.Range('E43').Value = '  +3.76%  '
$ws.Range('D44').Value = '''2.90'
$ws.Range('E44').Value = '  +21.01%  '
$ws.Range('E45').Value = '  +1.98%  '
$ws.Range('D46').Value = '1.351.15'
$ws.Range('E46').Value = '  +2.75%  '
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('E48').Value = '  +1.46%  '
$ws.Range('D49').Value = '''2.84'
$ws.Range('E49').Value = '  +3.42%  '
$ws.Range('E50').Value = '  +2.00%  '
$ws.Range('D51').Value = '''42.99'
$ws.Range('E51').Value = '  +0.92%  '
